# Revert BEPEfCT to default settings
$wb = $excel.ActiveWorkbook

$wsBEPEfCT = $wb.Worksheets.Item("BEPEfCT")

# Set the boolean lever cell B2 back to its default value (0)
$wsBEPEfCT.Range("B2").Value = 0

# Select cell B3 on the BEPEfCT sheet (matches the saved selection state)
$wsBEPEfCT.Range("B3").Select()

# Restore the "About" sheet as the active/selected tab
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
